$d = $word.ActiveDocument

# The document's Title-styled paragraph currently contains the run of text
# "Clay Freeman" followed by the _GoBack bookmark start/end. That name is
# redundant with the document's author metadata, so remove the run/text,
# leaving the (now empty) paragraph with its bookmark intact.
$d.Content.Find.Execute("Clay Freeman", $false, $false, $false, $false,
                         $false, $true, 1, $false, "", 2)
